# Insert 3 new rows at row 180 (pushing existing rows 180..277 down to 183..280)
# and populate them with the new data described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 180.
$ws.Range("A180:R182").EntireRow.Insert()

# Row 180
$ws.Range("A180").Value2 = 3
$ws.Range("B180").Value2 = "Femacal de La Calera"
$ws.Range("C180").Value2 = "Coquimbo"
$ws.Range("D180").Value2 = 44438
$ws.Range("E180").Value2 = 5
$ws.Range("F180").Value2 = 100112045
$ws.Range("G180").Value2 = "Zapallo"
$ws.Range("H180").Value2 = "Camote"
$ws.Range("I180").Value2 = "1a (guarda)"
$ws.Range("J180").Value2 = 260
$ws.Range("K180").Value2 = 750
$ws.Range("L180").Value2 = 800
$ws.Range("M180").Value2 = 771
$ws.Range("N180").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O180").Value2 = "Provincia de Talca"
$ws.Range("P180").Value2 = 771
$ws.Range("Q180").Value2 = 1
$ws.Range("R180").Value2 = "Hortaliza"

# Row 181
$ws.Range("A181").Value2 = 3
$ws.Range("B181").Value2 = "Femacal de La Calera"
$ws.Range("C181").Value2 = "Coquimbo"
$ws.Range("D181").Value2 = 44438
$ws.Range("E181").Value2 = 5
$ws.Range("F181").Value2 = 100112045
$ws.Range("G181").Value2 = "Zapallo"
$ws.Range("H181").Value2 = "Camote"
$ws.Range("I181").Value2 = "2a (guarda)"
$ws.Range("J181").Value2 = 80
$ws.Range("K181").Value2 = 600
$ws.Range("L181").Value2 = 600
$ws.Range("M181").Value2 = 600
$ws.Range("N181").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O181").Value2 = "Provincia de Talca"
$ws.Range("P181").Value2 = 600
$ws.Range("Q181").Value2 = 1
$ws.Range("R181").Value2 = "Hortaliza"

# Row 182
$ws.Range("A182").Value2 = 3
$ws.Range("B182").Value2 = "Femacal de La Calera"
$ws.Range("C182").Value2 = "Coquimbo"
$ws.Range("D182").Value2 = 44438
$ws.Range("E182").Value2 = 5
$ws.Range("F182").Value2 = 100112045
$ws.Range("G182").Value2 = "Zapallo"
$ws.Range("H182").Value2 = "Paine"
$ws.Range("I182").Value2 = "1a (guarda)"
$ws.Range("J182").Value2 = 120
$ws.Range("K182").Value2 = 600
$ws.Range("L182").Value2 = 600
$ws.Range("M182").Value2 = 600
$ws.Range("N182").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O182").Value2 = "Provincia de Talca"
$ws.Range("P182").Value2 = 600
$ws.Range("Q182").Value2 = 1
$ws.Range("R182").Value2 = "Hortaliza"
